$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the three runs that made up the Calvino/Viganò sentence into
#    a single run (same character formatting) and drop the gramStart /
#    gramEnd proofErr markers that used to wrap "e “".  A Find/Replace
#    over the whole (unchanged) sentence forces the host to rebuild the
#    run, which naturally coalesces it into one <w:r> and discards the
#    now-orphaned proofErr markers.
# ---------------------------------------------------------------------
$calvinoText = "Un esponente del neorealismo è sicuramente Italo Calvino, con i suoi romanzi “Il sentiero dei nidi di ragno” e “Ultimo viene il corvo”, entrambi ispirati al tema della resistenza. Altra esponente è Renata Viganò, col il romanzo “L’Agnese va a morire”, che narra di una donna anziana che entra a far parte della resistenza."

$d.Content.Find.Execute($calvinoText, $false, $false, $false, $false, $false, $true, 1, $false, $calvinoText, 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the trailing empty paragraph to
#    the very start of the document (right before the title run).
#    A throw-away placeholder character is used because collapsed
#    (zero-length) bookmark ranges get mis-anchored by the host; typing
#    a character first, bookmarking the one-character range and then
#    deleting that character again leaves a true zero-length bookmark
#    in the right spot.
# ---------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$docStart = $d.Paragraphs.Item(1).Range.Start
$placeholder = $d.Range($docStart, $docStart)
$placeholder.InsertBefore("X")
$markerRange = $d.Range($docStart, $docStart + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$cleanupRange = $d.Range($docStart, $docStart + 1)
$cleanupRange.Text = ""

# ---------------------------------------------------------------------
# 3. Shrink the page margins to 0.5" (720 twips) on every side.
# ---------------------------------------------------------------------
$pageSetup = $d.Sections.Item(1).PageSetup
$pageSetup.TopMargin = 36
$pageSetup.BottomMargin = 36
$pageSetup.LeftMargin = 36
$pageSetup.RightMargin = 36
